$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'230.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.306"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05542"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.377"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.475"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.078"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'7FTXTokenFTT"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.7744"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1371"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07423"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03152"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.02943"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09254"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.001667"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.251"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04777"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005896"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'17OneONEWorstin24h"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.006213"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005235"
$ws.Range("D20").Style = "Normal"
$ws.Range("B21").Value = "'UpBots"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'0.007485"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'20UpBotsUBXTBestin24h"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'BitKan"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.001063"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'21BitKanKAN"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'NitroEx"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'22NitroExNTX"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'LEO"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'3.946"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'23LEOLEO"
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'BTSEToken"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'2.196"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'24BTSETokenBTSE"
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'BitpandaEcosystemToken"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'0.3323"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'25BitpandaEcosystemTokenBEST"
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'ProBitToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'0.1248"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'26ProBitTokenPROB"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03957"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007157"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "'CEJI"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.003499"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'41CEJICEJI"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'BKEXToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.1038"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'42BKEXTokenBKK"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009820"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005441"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.7848"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.04021"
$ws.Range("D48").Style = "Normal"
